$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of merged cells in row 2 (year headers) that need to be unmerged,
# with the hidden right-hand cell populated with the same year value.
$yearPairs = @(
    @{Left="C2"; Right="D2"},
    @{Left="E2"; Right="F2"},
    @{Left="G2"; Right="H2"},
    @{Left="I2"; Right="J2"},
    @{Left="K2"; Right="L2"},
    @{Left="M2"; Right="N2"},
    @{Left="O2"; Right="P2"},
    @{Left="Q2"; Right="R2"},
    @{Left="S2"; Right="T2"},
    @{Left="U2"; Right="V2"}
)

foreach ($pair in $yearPairs) {
    $leftCell = $ws.Range($pair.Left)
    $value = $leftCell.Value2
    $leftCell.MergeArea.UnMerge()
    $ws.Range($pair.Right).Value = $value
}

# Merged note cells S42:V42, S43:V43, S44:V44 - unmerge and copy the note
# text into the newly revealed cells T, U, V.
$noteRows = @(42, 43, 44)

foreach ($row in $noteRows) {
    $srcAddr = "S" + $row
    $srcCell = $ws.Range($srcAddr)
    $value = $srcCell.Value2
    $srcCell.MergeArea.UnMerge()
    $ws.Range("T" + $row).Value = $value
    $ws.Range("U" + $row).Value = $value
    $ws.Range("V" + $row).Value = $value
}
